$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New headers for the two additional columns (F = planned start date,
#    G = planned finish date).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "计划开始时间"
$ws.Range("G1").Value = "计划完成时间"

# ---------------------------------------------------------------------------
# 2. Fill in the planned start / finish date serials for every data row
#    (rows 2-16).
# ---------------------------------------------------------------------------
$planStart = @(42399,42399,42399,42459,42459,42459,42459,42459,42459,42459,42459,42459,42459,42459,42459)
$planEnd   = @(42460,42460,42460,42734,42734,42734,42734,42734,42734,42734,42734,42734,42734,42734,42734)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $planStart[$i]
    $ws.Cells.Item($row, 7).Value = $planEnd[$i]
}

# Apply a date number format to F2 only, then fan that exact formatting out
# to the rest of F2:G16 via copy/paste-format so every cell shares a single
# style entry instead of each cell getting its own.
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()
$ws.Range("F2:G16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. The "计划级别" (plan level) column B used to be shown as a percentage;
#    switch it back to a plain general number so it reads 1 / 0.8 / 0.6 ...
# ---------------------------------------------------------------------------
$ws.Range("B2:B16").NumberFormat = "general"

# ---------------------------------------------------------------------------
# 4. Column widths for the new / resized columns.
# ---------------------------------------------------------------------------
$ws.Columns(2).ColumnWidth = 8.2857142857143   # B -> 9
$ws.Columns(4).ColumnWidth = 11.0              # D -> 11.75
$ws.Columns(6).ColumnWidth = 12.4285714285714  # F -> 13.125
$ws.Columns(7).ColumnWidth = 14.1428571428571  # G -> 14.875
$ws.Columns(5).ColumnWidth = 8.8571428571429   # E -> 9.5

# ---------------------------------------------------------------------------
# 5. Selection moves to G19, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("G19").Select()

# ---------------------------------------------------------------------------
# 6. Page setup: A4 paper, portrait orientation.
# ---------------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
